$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been recorded.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: Based on the conversation, the committee has decided to involve the rest of the PTA for input and has not come to a definitive decision about the movie to show on Friday. Therefore, I will call the no_decision function.`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The rights for `"Oppenheimer`" have been successfully acquired for the upcoming screening.`n"

$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been successfully recorded for showing `"Barbie`".`n"
$ws.Range("D5").Value = "Barbie_was_selected, "

$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday could not be finalized, leading to no selection being made.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"

$ws.Range("C8").Value = "MSG: None`n`nMSG: I have successfully recorded the decision to acquire the rights for both movies.`n"
$ws.Range("D8").Value = "both_movies, "

$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision about what movie to show on Friday has not been made.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday could not be made, so no selection will be made at this time.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The conversation has concluded without a decision about which movie to show on Friday.`n"
